$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of speaker id/variant pairs for rows 2-10 (re-exported
# without is_prefered flags and without Levenshtein-distance based matches)
$data = @(
    @("#philippyn", "Philippyn"),
    @("#lopes",      "Lopes"),
    @("#ferdinand",  "Ferdinand"),
    @("#philidia",   "Philidia"),
    @("#waardin",    "Waardin"),
    @("#alonce",     "Alonce"),
    @("#frederik",   "Frederik"),
    @("#elize",      "Elize"),
    @("#izabel",     "Izabel")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
    $ws.Cells.Item($row, 4).ClearContents()
    $row = $row + 1
}
